$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "2013 Verkaufszahlen (Detail)"
